# Update DateBase/orders/Fresh bloom Flowers_2025-10-14.xlsx
# - Append new order rows (32-41) to the "Orders" sheet
# - Refresh the concatenated "TotalNumber" digest on the "Summary" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New flower-order rows. Numbers are stored as text (matching the rest of
# the sheet), so force a text number format before assigning the values.
$newRows = @(
  @{ Row = 32; C = "479_绿灵草_lepidium_undefined_1bunch"; F = "5" },
  @{ Row = 33; C = "350_千层金红_Melaleuca bracteata`n（dyed red）_Melaleuca bracteata F.Muell._1bunch"; F = "10" },
  @{ Row = 34; C = "401_大飞燕白色_delphinium white_undefined_1bunch"; F = "10" },
  @{ Row = 35; C = "403_大飞燕浅蓝色_delphinium light blue_undefined_1bunch"; F = "10" },
  @{ Row = 36; C = "402_大飞燕深紫色_delphinium purple_undefined_1bunch"; F = "10" },
  @{ Row = 37; C = "305_彩星 宝蓝_Tinted Gypso blue_undefined_0.5kg"; F = "40" },
  @{ Row = 38; C = "301_彩星 红_Tinted Gypso Red_ gypsophila_0.5kg"; F = "10" },
  @{ Row = 39; C = "302_彩星 浅粉_Tinted Gypso light pink_undefined_0.5kg"; F = "4" },
  @{ Row = 40; C = "303_彩星 粉_Tinted Gypso pink_undefined_0.5kg"; F = "4" }
)

foreach ($r in $newRows) {
    $cCell = $ws.Cells.Item($r.Row, 3)   # column C = FlowerName
    $cCell.NumberFormat = "@"
    $cCell.Value = $r.C

    $fCell = $ws.Cells.Item($r.Row, 6)   # column F = Number
    $fCell.NumberFormat = "@"
    $fCell.Value = $r.F
}

# Trailing row 41 only carries a PackageID marker in column A.
$aCell = $ws.Cells.Item(41, 1)
$aCell.NumberFormat = "@"
$aCell.Value = "5"

# Refresh the Summary sheet's concatenated Number digest (column G, row 2)
# so it includes the newly added rows' Number values.
$ws2 = $wb.Worksheets.Item("Summary")
$g2 = $ws2.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "05881031020205558510312156555551710101510555101010104010440"
